$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3822
$ws.Range("I3").Value = 3957
$ws.Range("B4").Value = 1655
$ws.Range("D4").Value = 1927
$ws.Range("H4").Value = 1666
$ws.Range("I4").Value = 926
$ws.Range("I5").Value = 365
$ws.Range("I6").Value = 4436
$ws.Range("B7").Value = 23287
$ws.Range("D7").Value = 28117
$ws.Range("H7").Value = 25976
$ws.Range("I7").Value = 13506

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I4").Value = 52
$ws.Range("I7").Value = 426
$ws.Range("I8").Value = 818
$ws.Range("I10").Value = 91
$ws.Range("I14").Value = 70
$ws.Range("I15").Value = 160
$ws.Range("I17").Value = 16
$ws.Range("I19").Value = 368
$ws.Range("I20").Value = 330
$ws.Range("I22").Value = 38
$ws.Range("I23").Value = 132
$ws.Range("I29").Value = 874
$ws.Range("I31").Value = 130
$ws.Range("I33").Value = 612
$ws.Range("I36").Value = 190
$ws.Range("I37").Value = 439
$ws.Range("I42").Value = 473
$ws.Range("I43").Value = 118
$ws.Range("I51").Value = 130
$ws.Range("I54").Value = 307
$ws.Range("I55").Value = 147
$ws.Range("B63").Value = 361
$ws.Range("D63").Value = 318
$ws.Range("I63").Value = 48
$ws.Range("I65").Value = 299
$ws.Range("I67").Value = 522
$ws.Range("I68").Value = 44
$ws.Range("I70").Value = 25
$ws.Range("I72").Value = 50
$ws.Range("I76").Value = 204
$ws.Range("I78").Value = 198
$ws.Range("I79").Value = 364
$ws.Range("I83").Value = 272
$ws.Range("I85").Value = 610
$ws.Range("I90").Value = 169
$ws.Range("I91").Value = 167
$ws.Range("I95").Value = 215
$ws.Range("I96").Value = 149
$ws.Range("I98").Value = 89
$ws.Range("H99").Value = 438
$ws.Range("I99").Value = 251
$ws.Range("B101").Value = 23287
$ws.Range("D101").Value = 28117
$ws.Range("H101").Value = 25976
$ws.Range("I101").Value = 13506

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 247
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 610

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 253
$ws.Range("I3").Value = 228
$ws.Range("I6").Value = 266
$ws.Range("I7").Value = 818

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 147
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 426

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 140
$ws.Range("I3").Value = 140
$ws.Range("I7").Value = 439

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 90
$ws.Range("H4").Value = 25
$ws.Range("I5").Value = 7
$ws.Range("H7").Value = 438
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 187
$ws.Range("I6").Value = 172
$ws.Range("I7").Value = 522

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 97
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 299

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 95
$ws.Range("I7").Value = 272

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 78
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 145
$ws.Range("I3").Value = 223
$ws.Range("I6").Value = 193
$ws.Range("I7").Value = 612

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 60
$ws.Range("I6").Value = 154
$ws.Range("I7").Value = 307

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 258
$ws.Range("I3").Value = 300
$ws.Range("I6").Value = 239
$ws.Range("I7").Value = 874

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 138
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 368

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 50
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 129
$ws.Range("I3").Value = 163
$ws.Range("I6").Value = 127
$ws.Range("I7").Value = 473

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 50
$ws.Range("I4").Value = 26
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 147

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 35
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 52
$ws.Range("I3").Value = 58
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 109
$ws.Range("I3").Value = 113
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 364

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 92
$ws.Range("I3").Value = 101
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 330

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I3").Value = 8
$ws.Range("I7").Value = 16

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 60
$ws.Range("I3").Value = 59
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 89

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I2").Value = 16
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 52
